# Apply the commit's changes to the "总体功能列表.xlsx" workbook:
#  - update a handful of cells on "总体需求" (text replaced / cleared)
#  - add a new, mostly-empty sheet "第一个版本" after "总体需求"
#  - restore the selections on both sheets

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("总体需求")

# ---- cell content edits on 总体需求 --------------------------------------
# A12 used to read "按照板块筛选文章"; it now holds a small JS-ish snippet
$ws.Range("A12").Value = "{ count: Number }"

# The whole "本地图片/外源图片" design note block got cut - clear the cells
# that used to carry that text (keeping whatever cell formatting is there).
$ws.Range("B22").Value = ""
$ws.Range("B23").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("C21").Value = ""
$ws.Range("C23").Value = ""
$ws.Range("C24").Value = ""
$ws.Range("C25").Value = ""
$ws.Range("C26").Value = ""
$ws.Range("C29").Value = ""
$ws.Range("C30").Value = ""
$ws.Range("C31").Value = ""
$ws.Range("C32").Value = ""
$ws.Range("C33").Value = ""
$ws.Range("C34").Value = ""

# ---- add the new "第一个版本" sheet, right after 总体需求 -----------------
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "第一个版本"
$ws2.Columns.Item(1).ColumnWidth = 8.3
$ws2.Range("D25").Select()

# ---- restore selections / active sheet ------------------------------------
$ws.Activate()
$ws.Range("A12").Select()
